# Weekly update: insert a new price record for "Puerro" (Vega Modelo de Temuco)
# at row 112, pushing the existing rows 112-162 down to 113-163.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 112 (rows 112..162 shift down to 113..163).
$ws.Rows.Item(112).EntireRow.Insert()

# Populate the newly inserted row 112 with this week's new record.
$ws.Cells.Item(112, 1).Value = 10
$ws.Cells.Item(112, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(112, 3).Value = "La Araucanía"
$ws.Cells.Item(112, 4).Value = 44529
$ws.Cells.Item(112, 5).Value = 9
$ws.Cells.Item(112, 6).Value = 100112005
$ws.Cells.Item(112, 7).Value = "Puerro"
$ws.Cells.Item(112, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 10).Value = 65
$ws.Cells.Item(112, 11).Value = 8000
$ws.Cells.Item(112, 12).Value = 8000
$ws.Cells.Item(112, 13).Value = 8000
$ws.Cells.Item(112, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(112, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(112, 16).Value = 667
$ws.Cells.Item(112, 17).Value = 12
$ws.Cells.Item(112, 18).Value = "Hortaliza"

# Date column (D) uses the custom date/time number format applied to the rest
# of column D; make sure the new row matches it.
$ws.Cells.Item(112, 4).NumberFormat = $ws.Cells.Item(113, 4).NumberFormat
